$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17 (sheet ALC, diff line 1486)
$ws.Range("H17").Value = 1371316.1
$ws.Range("J17").Value = 1409925.6
$ws.Range("L17").Value = 4229776.800000001
$ws.Range("N17").Value = -4230112.800000001

# row 18 (sheet ALC, diff line 1538)
$ws.Range("H18").Value = 149
$ws.Range("I18").Value = 149
$ws.Range("K18").Value = 149
$ws.Range("M18").Value = 135

# row 33 (sheet ALC, diff line 2267)
$ws.Range("H33").Value = 87
$ws.Range("I33").Value = 87
$ws.Range("K33").Value = 87
$ws.Range("M33").Value = 142

# row 70 (sheet ALC, diff line 4134)
$ws.Range("H70").Value = 846.9
$ws.Range("J70").Value = 845
$ws.Range("L70").Value = 2535
$ws.Range("N70").Value = -3075

# row 73 (sheet ALC, diff line 4287)
$ws.Range("H73").Value = 846.9
$ws.Range("J73").Value = 845
$ws.Range("L73").Value = 2535
$ws.Range("N73").Value = -4407

# row 96 (sheet ALC, diff line 5426)
$ws.Range("H96").Value = 50000904
$ws.Range("I96").Value = 50000904
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 150002712
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -150001339

# row 100 (sheet ALC, diff line 5631)
$ws.Range("H100").Value = 142859650
$ws.Range("I100").Value = 333333980
$ws.Range("J100").Value = 3875
$ws.Range("K100").Value = 333333980
$ws.Range("L100").Value = 3875
$ws.Range("M100").Value = -333333439
$ws.Range("N100").Value = -4957

# row 116 (sheet ALC, diff line 6436)
$ws.Range("H116").Value = 3626.6316
$ws.Range("I116").Value = 1499.375
$ws.Range("K116").Value = 1499.375
$ws.Range("M116").Value = 1942.625

# row 125 (sheet ALC, diff line 6877)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# row 132 (sheet ALC, diff line 7223)
$ws.Range("H132").Value = 4046.6956
$ws.Range("I132").Value = 4010.1904
$ws.Range("J132").Value = 4430
$ws.Range("K132").Value = 12030.5712
$ws.Range("L132").Value = 13290
$ws.Range("M132").Value = -9500.5712
$ws.Range("N132").Value = -18350

# row 138 (sheet ALC, diff line 7526)
$ws.Range("H138").Value = 1362.8842
$ws.Range("I138").Value = 544.3111
$ws.Range("J138").Value = 2099.6
$ws.Range("K138").Value = 1632.9333
$ws.Range("L138").Value = 6298.799999999999
$ws.Range("M138").Value = 3507.0667
$ws.Range("N138").Value = -16578.8

$ws = $wb.Worksheets.Item("ARM")
# row 2 (sheet ARM, diff line 7819)
$ws.Range("H2").Value = 1296.4
$ws.Range("I2").Value = 1083.2
$ws.Range("K2").Value = 1083.2
$ws.Range("M2").Value = -970.2

# row 25 (sheet ARM, diff line 8934)
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4598

# row 97 (sheet ARM, diff line 12408)
$ws.Range("H97").Value = 862.7619
$ws.Range("I97").Value = 894.375
$ws.Range("J97").Value = 761.6
$ws.Range("K97").Value = 894.375
$ws.Range("L97").Value = 761.6
$ws.Range("M97").Value = -398.375
$ws.Range("N97").Value = -1753.6

# row 114 (sheet ARM, diff line 13235)
$ws.Range("H114").Value = 28419.8
$ws.Range("J114").Value = 28419.8
$ws.Range("L114").Value = 28419.8
$ws.Range("N114").Value = -37097.8

# row 116 (sheet ARM, diff line 13333)
$ws.Range("H116").Value = 1296.4
$ws.Range("I116").Value = 1083.2
$ws.Range("K116").Value = 1083.2
$ws.Range("M116").Value = 1210.8

# row 122 (sheet ARM, diff line 13627)
$ws.Range("H122").Value = 2606.9
$ws.Range("I122").Value = 1934
$ws.Range("K122").Value = 5802
$ws.Range("M122").Value = -3352

$ws = $wb.Worksheets.Item("BSM")
# row 3 (sheet BSM, diff line 14750)
$ws.Range("H3").Value = 1296.4
$ws.Range("I3").Value = 1083.2
$ws.Range("K3").Value = 1083.2
$ws.Range("M3").Value = -969.2

# row 94 (sheet BSM, diff line 19146)
$ws.Range("H94").Value = 1354.4839
$ws.Range("I94").Value = 1188.76
$ws.Range("J94").Value = 2045
$ws.Range("K94").Value = 1188.76
$ws.Range("L94").Value = 2045
$ws.Range("M94").Value = -737.76
$ws.Range("N94").Value = -2947

# row 99 (sheet BSM, diff line 19397)
$ws.Range("H99").Value = 1547.2307
$ws.Range("I99").Value = 1111.4
$ws.Range("K99").Value = 1111.4
$ws.Range("M99").Value = 386.5999999999999

# row 126 (sheet BSM, diff line 20699)
$ws.Range("H126").Value = 37195
$ws.Range("J126").Value = 37195
$ws.Range("L126").Value = 37195
$ws.Range("N126").Value = -47075

$ws = $wb.Worksheets.Item("CRP")
# row 31 (sheet CRP, diff line 22998)
$ws.Range("H31").Value = 12341.286
$ws.Range("I31").Value = 18834.217
$ws.Range("K31").Value = 18834.217
$ws.Range("M31").Value = -18539.217

# row 33 (sheet CRP, diff line 23099)
$ws.Range("H33").Value = 4666.6665
$ws.Range("I33").Value = 4666.6665
$ws.Range("K33").Value = 4666.6665
$ws.Range("M33").Value = -4287.6665

# row 34 (sheet CRP, diff line 23148)
$ws.Range("H34").Value = 12341.286
$ws.Range("I34").Value = 18834.217
$ws.Range("K34").Value = 18834.217
$ws.Range("M34").Value = -18632.217

# row 99 (sheet CRP, diff line 26306)
$ws.Range("H99").Value = 17245450
$ws.Range("I99").Value = 3376.5
$ws.Range("J99").Value = 38466464
$ws.Range("K99").Value = 3376.5
$ws.Range("L99").Value = 38466464
$ws.Range("M99").Value = -1878.5
$ws.Range("N99").Value = -38469460

# row 122 (sheet CRP, diff line 27418)
$ws.Range("H122").Value = 1045.9445
$ws.Range("I122").Value = 1033.3334
$ws.Range("K122").Value = 3100.0002
$ws.Range("M122").Value = -650.0001999999999

# row 126 (sheet CRP, diff line 27617)
$ws.Range("H126").Value = 17245450
$ws.Range("I126").Value = 3376.5
$ws.Range("J126").Value = 38466464
$ws.Range("K126").Value = 10129.5
$ws.Range("L126").Value = 115399392
$ws.Range("M126").Value = -7659.5
$ws.Range("N126").Value = -115404332

$ws = $wb.Worksheets.Item("CUL")
# row 5 (sheet CUL, diff line 28660)
$ws.Range("H5").Value = 681.5833
$ws.Range("J5").Value = 621.6667
$ws.Range("L5").Value = 1865.0001
$ws.Range("N5").Value = -2089.0001

# row 122 (sheet CUL, diff line 34543)
$ws.Range("H122").Value = 1353.6666
$ws.Range("I122").Value = 412.5
$ws.Range("J122").Value = 1575.1177
$ws.Range("K122").Value = 3712.5
$ws.Range("L122").Value = 14176.0593
$ws.Range("M122").Value = -1262.5
$ws.Range("N122").Value = -19076.0593

# row 131 (sheet CUL, diff line 34993)
$ws.Range("H131").Value = 781.47424
$ws.Range("J131").Value = 782.32294
$ws.Range("L131").Value = 2346.96882
$ws.Range("N131").Value = -12426.96882

# row 135 (sheet CUL, diff line 35201)
$ws.Range("H135").Value = 681.5833
$ws.Range("J135").Value = 621.6667
$ws.Range("L135").Value = 5595.0003
$ws.Range("N135").Value = -10665.0003

$ws = $wb.Worksheets.Item("GSM")
# row 74 (sheet GSM, diff line 39139)
$ws.Range("H74").Value = 47395
$ws.Range("J74").Value = 47395
$ws.Range("L74").Value = 47395
$ws.Range("N74").Value = -49267

# row 77 (sheet GSM, diff line 39280)
$ws.Range("H77").Value = 47395
$ws.Range("J77").Value = 47395
$ws.Range("L77").Value = 142185
$ws.Range("N77").Value = -151545

# row 102 (sheet GSM, diff line 40484)
$ws.Range("H102").Value = 20835340
$ws.Range("I102").Value = 25002226
$ws.Range("K102").Value = 25002226
$ws.Range("M102").Value = -25000604

# row 103 (sheet GSM, diff line 40536)
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

# row 111 (sheet GSM, diff line 40919)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0

# row 122 (sheet GSM, diff line 41443)
$ws.Range("H122").Value = 51283012
$ws.Range("I122").Value = 23810704
$ws.Range("J122").Value = 83334040
$ws.Range("K122").Value = 71432112
$ws.Range("L122").Value = 250002120
$ws.Range("M122").Value = -71429662
$ws.Range("N122").Value = -250007020

# row 126 (sheet GSM, diff line 41642)
$ws.Range("H126").Value = 4689.6206
$ws.Range("I126").Value = 3605.2632
$ws.Range("J126").Value = 6749.9
$ws.Range("K126").Value = 10815.7896
$ws.Range("L126").Value = 20249.7
$ws.Range("M126").Value = -8345.7896
$ws.Range("N126").Value = -25189.7

$ws = $wb.Worksheets.Item("LTW")
# row 25 (sheet LTW, diff line 43614)
$ws.Range("H25").Value = 2960
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 3440
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 3440
$ws.Range("M25").Value = -1770
$ws.Range("N25").Value = -3900

# row 40 (sheet LTW, diff line 44364)
$ws.Range("H40").Value = 3958.7778
$ws.Range("I40").Value = 2479.3333
$ws.Range("K40").Value = 2479.3333
$ws.Range("M40").Value = -2343.3333

# row 68 (sheet LTW, diff line 45736)
$ws.Range("H68").Value = 2753.8572
$ws.Range("J68").Value = 2469.25
$ws.Range("L68").Value = 2469.25
$ws.Range("N68").Value = -3967.25

# row 71 (sheet LTW, diff line 45883)
$ws.Range("H71").Value = 2753.8572
$ws.Range("J71").Value = 2469.25
$ws.Range("L71").Value = 12346.25
$ws.Range("N71").Value = -19834.25

# row 93 (sheet LTW, diff line 46928)
$ws.Range("H93").Value = 1236.7037
$ws.Range("I93").Value = 1374.8096
$ws.Range("J93").Value = 753.3333
$ws.Range("K93").Value = 1374.8096
$ws.Range("L93").Value = 753.3333
$ws.Range("M93").Value = -126.8096
$ws.Range("N93").Value = -3249.3333

# row 122 (sheet LTW, diff line 48322)
$ws.Range("H122").Value = 1785568
$ws.Range("I122").Value = 1963354.8
$ws.Range("K122").Value = 5890064.4
$ws.Range("M122").Value = -5887614.4

# row 136 (sheet LTW, diff line 49008)
$ws.Range("H136").Value = 43000.582
$ws.Range("I136").Value = 72714.57000000001
$ws.Range("K136").Value = 218143.71
$ws.Range("M136").Value = -215593.71

$ws = $wb.Worksheets.Item("WVR")
# row 24 (sheet WVR, diff line 50468)
$ws.Range("H24").Value = 20000000
$ws.Range("I24").Value = 20000000
$ws.Range("K24").Value = 20000000
$ws.Range("M24").Value = -19999770

# row 31 (sheet WVR, diff line 50799)
$ws.Range("H31").Value = 6750
$ws.Range("I31").Value = 3500
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 3500
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -3152
$ws.Range("N31").Value = -10696

# row 100 (sheet WVR, diff line 54150)
$ws.Range("H100").Value = 163.28572
$ws.Range("I100").Value = 192
$ws.Range("J100").Value = 91.5
$ws.Range("K100").Value = 384
$ws.Range("L100").Value = 183
$ws.Range("M100").Value = 157
$ws.Range("N100").Value = -1265

# row 122 (sheet WVR, diff line 55204)
$ws.Range("H122").Value = 1651
$ws.Range("I122").Value = 1649.76
$ws.Range("K122").Value = 4949.28
$ws.Range("M122").Value = -2499.28

# row 136 (sheet WVR, diff line 55890)
$ws.Range("H136").Value = 40002410
$ws.Range("I136").Value = 76925590
$ws.Range("K136").Value = 230776770
$ws.Range("M136").Value = -230774220
